$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk at old line 2607
$ws.Range("H40").Value = 2235.3333
$ws.Range("I40").Value = 2940.5833
$ws.Range("J40").Value = 1295
$ws.Range("K40").Value = 2940.5833
$ws.Range("L40").Value = 1295
$ws.Range("M40").Value = -2765.5833
$ws.Range("N40").Value = -1645
# hunk at old line 2763
$ws.Range("H43").Value = 5634.357
$ws.Range("I43").Value = 3474.25
$ws.Range("K43").Value = 3474.25
$ws.Range("M43").Value = -3405.25
# hunk at old line 3715
$ws.Range("H62").Value = 19234590
$ws.Range("I62").Value = 38463570
$ws.Range("J62").Value = 5612.231
$ws.Range("K62").Value = 38463570
$ws.Range("L62").Value = 5612.231
$ws.Range("M62").Value = -38462946
$ws.Range("N62").Value = -6860.231
# hunk at old line 3865
$ws.Range("H65").Value = 19234590
$ws.Range("I65").Value = 38463570
$ws.Range("J65").Value = 5612.231
$ws.Range("K65").Value = 192317850
$ws.Range("L65").Value = 28061.155
$ws.Range("M65").Value = -192314730
$ws.Range("N65").Value = -34301.155
# hunk at old line 4309
$ws.Range("H74").Value = 3918.6667
$ws.Range("I74").Value = 3920
$ws.Range("J74").Value = 3917.3333
$ws.Range("K74").Value = 3920
$ws.Range("L74").Value = 3917.3333
$ws.Range("M74").Value = -2984
$ws.Range("N74").Value = -5789.3333
# hunk at old line 4410
$ws.Range("H76").Value = 5006.8335
$ws.Range("I76").Value = 4561.857
$ws.Range("K76").Value = 4561.857
$ws.Range("M76").Value = -4246.857
# hunk at old line 4462
$ws.Range("H77").Value = 3918.6667
$ws.Range("I77").Value = 3920
$ws.Range("J77").Value = 3917.3333
$ws.Range("K77").Value = 19600
$ws.Range("L77").Value = 19586.6665
$ws.Range("M77").Value = -14920
$ws.Range("N77").Value = -28946.6665
# hunk at old line 4563
$ws.Range("H79").Value = 5006.8335
$ws.Range("I79").Value = 4561.857
$ws.Range("K79").Value = 4561.857
$ws.Range("M79").Value = -3469.857
# hunk at old line 5916
$ws.Range("H106").Value = 16430.223
$ws.Range("I106").Value = 1995.5
$ws.Range("J106").Value = 45299.668
$ws.Range("K106").Value = 1995.5
$ws.Range("L106").Value = 45299.668
$ws.Range("M106").Value = -1364.5
$ws.Range("N106").Value = -46561.668
# hunk at old line 7193
$ws.Range("H132").Value = 5103.467
$ws.Range("I132").Value = 2849.125
$ws.Range("J132").Value = 7679.857
$ws.Range("K132").Value = 8547.375
$ws.Range("L132").Value = 23039.571
$ws.Range("M132").Value = -6017.375
$ws.Range("N132").Value = -28099.571

$ws = $wb.Worksheets.Item("ARM")
# hunk at old line 9893
$ws.Range("H45").Value = 1547.4445
$ws.Range("I45").Value = 1418.1428
$ws.Range("K45").Value = 1418.1428
$ws.Range("M45").Value = -1041.1428
# hunk at old line 10665
$ws.Range("H61").Value = 4123.2173
$ws.Range("I61").Value = 1946.3334
$ws.Range("K61").Value = 1946.3334
$ws.Range("M61").Value = -1734.3334
# hunk at old line 10766
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
# hunk at old line 10913
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
# hunk at old line 12644
$ws.Range("H102").Value = 2883
$ws.Range("I102").Value = 2580.5715
$ws.Range("K102").Value = 2580.5715
$ws.Range("M102").Value = -958.5715
# hunk at old line 14292
$ws.Range("H136").Value = 4123.2173
$ws.Range("I136").Value = 1946.3334
$ws.Range("K136").Value = 5839.0002
$ws.Range("M136").Value = -3289.0002

$ws = $wb.Worksheets.Item("BSM")
# hunk at old line 19388
$ws.Range("H99").Value = 3813.9644
$ws.Range("I99").Value = 3950.9583
$ws.Range("K99").Value = 3950.9583
$ws.Range("M99").Value = -2452.9583
# hunk at old line 19679
$ws.Range("H105").Value = 2090.087
$ws.Range("I105").Value = 2126
$ws.Range("K105").Value = 2126
$ws.Range("M105").Value = -379
# hunk at old line 19780
$ws.Range("H107").Value = 754.725
$ws.Range("I107").Value = 652.85187
$ws.Range("J107").Value = 966.3077
$ws.Range("K107").Value = 652.85187
$ws.Range("L107").Value = 966.3077
$ws.Range("M107").Value = 1267.14813
$ws.Range("N107").Value = -4806.3077

$ws = $wb.Worksheets.Item("CRP")
# hunk at old line 22539
$ws.Range("H22").Value = 302.05713
$ws.Range("I22").Value = 304.1875
$ws.Range("J22").Value = 279.33334
$ws.Range("K22").Value = 304.1875
$ws.Range("L22").Value = 279.33334
$ws.Range("M22").Value = 45.8125
$ws.Range("N22").Value = -979.33334
# hunk at old line 22995
$ws.Range("H31").Value = 372484.06
$ws.Range("I31").Value = 667698.9
$ws.Range("K31").Value = 667698.9
$ws.Range("M31").Value = -667403.9
# hunk at old line 23151
$ws.Range("H34").Value = 372484.06
$ws.Range("I34").Value = 667698.9
$ws.Range("K34").Value = 667698.9
$ws.Range("M34").Value = -667496.9
# hunk at old line 26315
$ws.Range("H99").Value = 951091.4399999999
$ws.Range("I99").Value = 3343337.2
$ws.Range("J99").Value = 53999.25
$ws.Range("K99").Value = 3343337.2
$ws.Range("L99").Value = 53999.25
$ws.Range("M99").Value = -3341839.2
$ws.Range("N99").Value = -56995.25
# hunk at old line 27626
$ws.Range("H126").Value = 951091.4399999999
$ws.Range("I126").Value = 3343337.2
$ws.Range("J126").Value = 53999.25
$ws.Range("K126").Value = 10030011.6
$ws.Range("L126").Value = 161997.75
$ws.Range("M126").Value = -10027541.6
$ws.Range("N126").Value = -166937.75

$ws = $wb.Worksheets.Item("GSM")
# hunk at old line 38901
$ws.Range("H70").Value = 8425.700000000001
$ws.Range("I70").Value = 8327
$ws.Range("K70").Value = 8327
$ws.Range("M70").Value = -8057
# hunk at old line 39045
$ws.Range("H73").Value = 8425.700000000001
$ws.Range("I73").Value = 8327
$ws.Range("K73").Value = 8327
$ws.Range("M73").Value = -7391
# hunk at old line 39373
$ws.Range("H80").Value = 4227.148
$ws.Range("I80").Value = 3010.25
$ws.Range("J80").Value = 4739.5264
$ws.Range("K80").Value = 3010.25
$ws.Range("L80").Value = 4739.5264
$ws.Range("M80").Value = -2012.25
$ws.Range("N80").Value = -6735.5264
# hunk at old line 39520
$ws.Range("H83").Value = 4227.148
$ws.Range("I83").Value = 3010.25
$ws.Range("J83").Value = 4739.5264
$ws.Range("K83").Value = 15051.25
$ws.Range("L83").Value = 23697.632
$ws.Range("M83").Value = -10059.25
$ws.Range("N83").Value = -33681.632
# hunk at old line 40424
$ws.Range("H102").Value = 3575
$ws.Range("I102").Value = 3695.4666
$ws.Range("K102").Value = 3695.4666
$ws.Range("M102").Value = -2073.4666
# hunk at old line 40954
$ws.Range("H113").Value = 14534.615
$ws.Range("J113").Value = 10987
$ws.Range("L113").Value = 10987
$ws.Range("N113").Value = -15327
# hunk at old line 41383
$ws.Range("H122").Value = 5025.75
$ws.Range("I122").Value = 5595.8
$ws.Range("J122").Value = 3315.6
$ws.Range("K122").Value = 16787.4
$ws.Range("L122").Value = 9946.799999999999
$ws.Range("M122").Value = -14337.4
$ws.Range("N122").Value = -14846.8

$ws = $wb.Worksheets.Item("LTW")
# hunk at old line 44475
$ws.Range("H43").Value = 16757.125
$ws.Range("J43").Value = 12905.667
$ws.Range("L43").Value = 12905.667
$ws.Range("N43").Value = -13291.667
# hunk at old line 45718
$ws.Range("H68").Value = 3790.8635
$ws.Range("I68").Value = 4312.4375
$ws.Range("K68").Value = 4312.4375
$ws.Range("M68").Value = -3563.4375
# hunk at old line 45862
$ws.Range("H71").Value = 3790.8635
$ws.Range("I71").Value = 4312.4375
$ws.Range("K71").Value = 21562.1875
$ws.Range("M71").Value = -17818.1875
# hunk at old line 46343
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
# hunk at old line 46392
$ws.Range("H82").Value = 1914
$ws.Range("I82").Value = 1905.2941
$ws.Range("J82").Value = 1988
$ws.Range("K82").Value = 1905.2941
$ws.Range("L82").Value = 1988
$ws.Range("M82").Value = -1544.2941
$ws.Range("N82").Value = -2710
# hunk at old line 46493
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
# hunk at old line 46542
$ws.Range("H85").Value = 1914
$ws.Range("I85").Value = 1905.2941
$ws.Range("J85").Value = 1988
$ws.Range("K85").Value = 1905.2941
$ws.Range("L85").Value = 1988
$ws.Range("M85").Value = -657.2941000000001
$ws.Range("N85").Value = -4484
# hunk at old line 48806
$ws.Range("H132").Value = 2106.543
$ws.Range("I132").Value = 1949.16
$ws.Range("K132").Value = 5847.48
$ws.Range("M132").Value = -3317.48

$ws = $wb.Worksheets.Item("WVR")
# hunk at old line 55216
$ws.Range("H122").Value = 24262.23
$ws.Range("I122").Value = 2613.3
$ws.Range("K122").Value = 7839.900000000001
$ws.Range("M122").Value = -5389.900000000001
# hunk at old line 55415
$ws.Range("H126").Value = 4008.111
$ws.Range("I126").Value = 3996.7144
$ws.Range("J126").Value = 4048
$ws.Range("K126").Value = 11990.1432
$ws.Range("L126").Value = 12144
$ws.Range("M126").Value = -9520.143199999999
$ws.Range("N126").Value = -17084
